$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '29.497.28', '  -1.43%  ', 0),
    @('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.851.74', '  -0.71%  ', 0),
    @('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9993', '  -0.14%  ', 1),
    @('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '243.14', '  -1.27%  ', 1),
    @('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.6379', '  -0.40%  ', 1),
    @('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.000', '  -0.04%  ', 1),
    @('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2983', '  -0.69%  ', 1),
    @('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07460', '  -0.34%  ', 1),
    @('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '24.22', '  -0.68%  ', 1),
    @('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07625', '  -0.71%  ', 1),
    @('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.848.17', '  -1.12%  ', 0),
    @('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.021', '  -1.15%  ', 1),
    @('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6839', '  -0.96%  ', 1),
    @('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '83.49', '  -0.99%  ', 1),
    @('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000009543', '  +0.81%  ', 1),
    @('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.133', '  +0.42%  ', 1),
    @('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '29.511.11', '  -1.36%  ', 0),
    @('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.071.59', '  -2.66%  ', 0),
    @('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '235.12', '  -2.75%  ', 1),
    @('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '12.53', '  -1.34%  ', 1),
    @('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.000', '  -0.06%  ', 1),
    @('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.663', '  +2.32%  ', 1),
    @('BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.001', '  -0.08%  ', 1),
    @('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '156.93', '  -1.80%  ', 1),
    @('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1405', '  -1.46%  ', 1),
    @('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.466', '  -1.55%  ', 1),
    @('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.75', '  -1.71%  ', 1),
    @('PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.486', '  -1.50%  ', 1),
    @('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05992', '  -2.74%  ', 1),
    @('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.255', '  -1.62%  ', 1),
    @('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.116', '  -1.31%  ', 1),
    @('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.068', '  -1.70%  ', 1),
    @('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.865', '  -0.31%  ', 1),
    @('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.177', '  +1.02%  ', 1),
    @('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7177', '  -2.58%  ', 1),
    @('HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.599', '  -0.31%  ', 1),
    @('MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.802', '  -2.39%  ', 1),
    @('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01773', '  -1.84%  ', 1),
    @('Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.195.83', '  -2.32%  ', 0),
    @('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.9087', '  -2.08%  ', 1),
    @('FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.160', '  -2.15%  ', 1),
    @('PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.0000', '  -0.26%  ', 1),
    @('RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '2.005.65', '  -1.65%  ', 0),
    @('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '101.76', '  -0.57%  ', 1),
    @('Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '66.25', '  -0.85%  ', 1),
    @('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.291', '  +8.24%  ', 1),
    @('BabyDogeCoin', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge', '0.00000000121', '  -2.58%  ', 1),
    @('TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4030', '  -1.76%  ', 1),
    @('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.117', '  -2.92%  ', 1),
    @('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.657', '  +1.42%  ', 1),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $dCell = $ws.Cells.Item($r, 4)
    if ($row[4] -eq 1) {
        $dCell.NumberFormat = "@"
    }
    $dCell.Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}
